# Apply BOM update: update the C3 component line to the new part,
# remove the old C4 "22uF 16V Aluminum Electrolytic Capacitor" row
# (shifting later rows up), and move the active selection to A8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 3 (C3 part) description / footprint / LCSC part number.
$ws.Range("A3").Value = "10uF 50V Multilayer Ceramic Capacitor"
$ws.Range("C3").Value = "0805"
$ws.Range("D3").Value = "C2932476"

# Delete the old row 4 (22uF 16V Aluminum Electrolytic Capacitor / C4)
# entirely; rows below shift up.
$ws.Rows("4").Delete()

# Update the saved selection/active cell to A8, matching the new layout.
$ws.Range("A8").Select()
